$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.260.88"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").Value = "1.920.94"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -1.50%  "
$ws.Range("D5").Value = "'316.76"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("D7").Value = "'0.4852"
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("D9").Value = "'0.07420"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "'0.9515"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").Value = "'20.97"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "'0.07807"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "1.923.29"
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").Value = "'5.558"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").Value = "'6.670"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "'92.23"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "'0.000008892"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "28.258.55"
$ws.Range("E20").Value = "  +3.34%  "
$ws.Range("D21").Value = "'15.06"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "'5.179"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "2.158.35"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").Value = "'10.95"
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").Value = "'156.31"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "'18.67"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "'2.118"
$ws.Range("E28").Value = "  +5.16%  "
$ws.Range("D29").Value = "'117.45"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").Value = "'5.031"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").Value = "'0.08912"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'3.366"
$ws.Range("D33").Value = "'1.253"
$ws.Range("E33").Value = "  +4.52%  "
$ws.Range("D34").Value = "'0.7778"
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("D35").Value = "'4.690"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").Value = "'2.785"
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "'1.132"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5603"
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.05379"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").Value = "'7.136"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'8.593"
$ws.Range("E43").Value = "  +2.34%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'10.82"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").Value = "'0.4926"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").Value = "'107.75"
$ws.Range("E47").Value = "  +4.44%  "
$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").Value = "'69.54"
$ws.Range("D51").Value = "'0.06150"
$ws.Range("E51").Value = "  +1.00%  "
